$wb = $excel.ActiveWorkbook

# --- Overview sheet: mark the handback-transform-failed status for the
#     9f2fe086 file (row 3) in both language status columns.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

# --- zh-cn sheet: update Status for the 9f2fe086 row (row 3) and populate the
#     Error Detail column (P) for that same row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("P3").Value = "Handback file name: pvzysdmv.4wn is different with handoff file name: 9f2fe086-89f0-4bea-a8a9-1c145fdf0e44.3d54703557453feaf74dbb3f4628adf4063cf908.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1666666666667

# --- de-de sheet: update Status for the 9f2fe086 row (row 3) and populate the
#     Error Detail column (P) for that same row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("P3").Value = "Handback file name: pvzysdmv.4wn is different with handoff file name: 9f2fe086-89f0-4bea-a8a9-1c145fdf0e44.3d54703557453feaf74dbb3f4628adf4063cf908.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1666666666667
